# Datos_Problema-11.xlsx : "Modifica datos del problema 11"
#
# Summary of the edit:
#  1. Rename sheet "Costos-Capacidades" -> "Costos" and drop its
#     "Capacidad" column (column D), reusing the vacated "Samborondon"
#     cost row as a new "Ibarra" entry so the city list stays in
#     alphabetical order.
#  2. "Otros-parametros" becomes the active/selected sheet; the
#     "Presupuesto" parameter is renamed "Presupuesto 01" (budget raised
#     to 2'500.000) and two new budget tiers are appended
#     ("Presupuesto 02" = 5'000.000, "Presupuesto 03" = 10'000.000);
#     the "Area Fija" parameter is raised from 200 to 220.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Costos-Capacidades -> Costos : drop the Capacidad column and
#    replace "Samborondon" with "Ibarra" (keeping the alphabetical order
#    of the city list, and reusing Samborondon's old cost value).
# ---------------------------------------------------------------------
$wsCostos = $wb.Worksheets.Item("Costos-Capacidades")

# Remove column D ("Capacidad") entirely.
$wsCostos.Columns.Item(4).Delete()

# Shift the Machala..Salinas rows down by one and insert "Ibarra" (cost
# 162.46, previously Samborondon's cost) in its alphabetically-correct
# spot just before Machala.
$wsCostos.Range("B6").Value = "Ibarra"
$wsCostos.Range("C6").Value = 162.46

$wsCostos.Range("B7").Value = "Machala"
$wsCostos.Range("C7").Value = 557.34

$wsCostos.Range("B8").Value = "Manta"
$wsCostos.Range("C8").Value = 515.08

$wsCostos.Range("B9").Value = "Quevedo"
$wsCostos.Range("C9").Value = 286.82

$wsCostos.Range("B10").Value = "Quito"
$wsCostos.Range("C10").Value = 244.48

$wsCostos.Range("B11").Value = "Riobamba"
$wsCostos.Range("C11").Value = 497.22

$wsCostos.Range("B12").Value = "Salinas"
$wsCostos.Range("C12").Value = 572.32

# Row 13 (Santo Domingo de Los Colorados / 249.07) is unchanged.

$wsCostos.Name = "Costos"

# ---------------------------------------------------------------------
# 2. Otros-parametros: update budget data and add two new budget rows.
# ---------------------------------------------------------------------
$wsOtros = $wb.Worksheets.Item("Otros-parametros")

# "Area Fija" value raised from 200 to 220.
$wsOtros.Range("B2").Value = 220

# "Presupuesto" -> "Presupuesto 01", value raised from 2'400.000 to 2'500.000.
$wsOtros.Range("A4").Value = "Presupuesto 01"
$wsOtros.Range("B4").Value = 2500000
$wsOtros.Range("C4").Value = "USD"

# New rows 5 and 6, copying row 4's formatting first so the new budget
# tiers look the same as the existing one.
$wsOtros.Range("A4:C4").Copy()
$wsOtros.Range("A5:C5").PasteSpecial(-4122)
$wsOtros.Range("A6:C6").PasteSpecial(-4122)

$wsOtros.Range("A5").Value = "Presupuesto 02"
$wsOtros.Range("B5").Value = 5000000
$wsOtros.Range("C5").Value = "USD"

$wsOtros.Range("A6").Value = "Presupuesto 03"
$wsOtros.Range("B6").Value = 10000000
$wsOtros.Range("C6").Value = "USD"

# "Otros-parametros" becomes the selected / active sheet.
$wsOtros.Activate()
